$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1), columns B..G (existing cells, keep styling) ---
$ws.Cells.Item(1, 2).Value = "Algorithm"
$ws.Cells.Item(1, 3).Value = "One Year Alt mean"
$ws.Cells.Item(1, 4).Value = "One Year Alt std"
$ws.Cells.Item(1, 5).Value = "Two Year Alt mean"
$ws.Cells.Item(1, 6).Value = "Two Year Alt std"
$ws.Cells.Item(1, 7).Value = "Three Year Alt mean"
$ws.Cells.Item(1, 8).Value = "Three Year Alt std"
$ws.Cells.Item(1, 9).Value = "Five Year Alt mean"
$ws.Cells.Item(1, 10).Value = "Five Year Alt std"
$ws.Cells.Item(1, 11).Value = "Ten Year Alt mean"
$ws.Cells.Item(1, 12).Value = "Ten Year Alt std"

# --- Add new header cells H1:L1 and copy the bold/border style from G1 ---
$ws.Range("G1").Copy()
$ws.Range("H1:L1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update algorithm name column (B) rows 2..8 ---
$ws.Cells.Item(2, 2).Value = "LR"
$ws.Cells.Item(3, 2).Value = "LDA"
$ws.Cells.Item(4, 2).Value = "KNN"
$ws.Cells.Item(5, 2).Value = "DTREE"
$ws.Cells.Item(6, 2).Value = "RTREE"
$ws.Cells.Item(7, 2).Value = "XTREE"
$ws.Cells.Item(8, 2).Value = "SVM"

# --- Update numeric data columns C..L rows 2..8 ---
$ws.Cells.Item(2, 3).Value = 0.8530708001675743
$ws.Cells.Item(2, 4).Value = 0.01961582804224883
$ws.Cells.Item(2, 5).Value = 0.848234969663541
$ws.Cells.Item(2, 6).Value = 0.02792517416860756
$ws.Cells.Item(2, 7).Value = 0.8367845370092898
$ws.Cells.Item(2, 8).Value = 0.02967221244273211
$ws.Cells.Item(2, 9).Value = 0.8304869186046512
$ws.Cells.Item(2, 10).Value = 0.05196904854839265
$ws.Cells.Item(2, 11).Value = 0.8348690991548133
$ws.Cells.Item(2, 12).Value = 0.04514654158089394

$ws.Cells.Item(3, 3).Value = 0.8362589023879347
$ws.Cells.Item(3, 4).Value = 0.01882206094725105
$ws.Cells.Item(3, 5).Value = 0.8224857510571797
$ws.Cells.Item(3, 6).Value = 0.03228151947333432
$ws.Cells.Item(3, 7).Value = 0.8212466287084206
$ws.Cells.Item(3, 8).Value = 0.0279470351565709
$ws.Cells.Item(3, 9).Value = 0.8079457364341085
$ws.Cells.Item(3, 10).Value = 0.05127120857689344
$ws.Cells.Item(3, 11).Value = 0.8206967635539065
$ws.Cells.Item(3, 12).Value = 0.0441683389691347

$ws.Cells.Item(4, 3).Value = 0.8886761625471303
$ws.Cells.Item(4, 4).Value = 0.02688024809305695
$ws.Cells.Item(4, 5).Value = 0.9044447508733224
$ws.Cells.Item(4, 6).Value = 0.02096313643083356
$ws.Cells.Item(4, 7).Value = 0.8961192688043151
$ws.Cells.Item(4, 8).Value = 0.01376519663223573
$ws.Cells.Item(4, 9).Value = 0.8958454457364342
$ws.Cells.Item(4, 10).Value = 0.02440718051138728
$ws.Cells.Item(4, 11).Value = 0.883395176252319
$ws.Cells.Item(4, 12).Value = 0.04904054332401061

$ws.Cells.Item(5, 3).Value = 0.7935860913280268
$ws.Cells.Item(5, 4).Value = 0.03523876775268898
$ws.Cells.Item(5, 5).Value = 0.7947049089906233
$ws.Cells.Item(5, 6).Value = 0.0282916357251466
$ws.Cells.Item(5, 7).Value = 0.7851613225452002
$ws.Cells.Item(5, 8).Value = 0.01676397973977825
$ws.Cells.Item(5, 9).Value = 0.7745881782945736
$ws.Cells.Item(5, 10).Value = 0.04043370704925065
$ws.Cells.Item(5, 11).Value = 0.7710265924551639
$ws.Cells.Item(5, 12).Value = 0.04612264472018746

$ws.Cells.Item(6, 3).Value = 0.860192710515291
$ws.Cells.Item(6, 4).Value = 0.02086526908164362
$ws.Cells.Item(6, 5).Value = 0.8523211987497701
$ws.Cells.Item(6, 6).Value = 0.03396143598093248
$ws.Cells.Item(6, 7).Value = 0.8424483068624513
$ws.Cells.Item(6, 8).Value = 0.02763488898733027
$ws.Cells.Item(6, 9).Value = 0.8320494186046512
$ws.Cells.Item(6, 10).Value = 0.04656977052918741
$ws.Cells.Item(6, 11).Value = 0.8277983920841063
$ws.Cells.Item(6, 12).Value = 0.04504070379367365

$ws.Cells.Item(7, 3).Value = 0.8453121072475911
$ws.Cells.Item(7, 4).Value = 0.01901348060213343
$ws.Cells.Item(7, 5).Value = 0.8346938775510203
$ws.Cells.Item(7, 6).Value = 0.04349618185646725
$ws.Cells.Item(7, 7).Value = 0.8290130856058336
$ws.Cells.Item(7, 8).Value = 0.03202015374506993
$ws.Cells.Item(7, 9).Value = 0.8453064437984497
$ws.Cells.Item(7, 10).Value = 0.04199537190734103
$ws.Cells.Item(7, 11).Value = 0.8561224489795919
$ws.Cells.Item(7, 12).Value = 0.04245767616234759

$ws.Cells.Item(8, 3).Value = 0.8795978215333055
$ws.Cells.Item(8, 4).Value = 0.01907669610893491
$ws.Cells.Item(8, 5).Value = 0.8760250045964332
$ws.Cells.Item(8, 6).Value = 0.02338973772457094
$ws.Cells.Item(8, 7).Value = 0.8820197782439317
$ws.Cells.Item(8, 8).Value = 0.03317021991192351
$ws.Cells.Item(8, 9).Value = 0.8841872577519381
$ws.Cells.Item(8, 10).Value = 0.02914437809416581
$ws.Cells.Item(8, 11).Value = 0.8662956091527519
$ws.Cells.Item(8, 12).Value = 0.03774202519599633

# --- Remove old row 9 (NB row no longer present in the new data) ---
$ws.Rows("9:9").Delete()